$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''257.27'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''0.97%'
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''27.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''-3.61%'
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''4.813'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''-10.26%'
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''0.05955'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''1.79%'
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''6.662'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''-0.67%'
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.8703'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''0.11%'
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.9593'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''4.28%'
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.1409'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''-0.79%'
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.03739'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''8.02%'
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.07176'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''0.66%'
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03179'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''-0.05%'
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09247'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''-0.02%'
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001554'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''-0.76%'
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = '''0.0006080'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''0.44%'
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.006067'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''4.36%'
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.480'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''-0.56%'
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '''3.190'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''-1.30%'
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = '''2.218'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''-0.17%'
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '''0.3133'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''-1.50%'
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = '''-0.75%'
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''3.809'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''7.63%'
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''0.04222'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''1.21%'
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '''0.02%'
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''0.001224'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''-0.39%'
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''0.004499'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''-10.57%'
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''0.0001690'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''40.78%'
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''-23.01%'
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = '''0.03839'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''0.22%'
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.006150'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''18.46%'
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''0.1100'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''-0.12%'
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''-4.57%'
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''-3.58%'
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.00005500'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''5.29%'
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''-0.05%'
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.08852'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''1.05%'
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''0.002366'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''9.61%'
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''-0.05%'
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''-0.05%'
$ws.Range("E50").Style = "Normal"
